$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.525.41"
$ws.Range("E2").Value = "  +5.59%  "

$ws.Range("D3").Value = "1.722.43"
$ws.Range("E3").Value = "  +4.36%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "'225.91"
$ws.Range("E5").Value = "  +3.62%  "

$ws.Range("D6").Value = "'0.5366"
$ws.Range("E6").Value = "  +3.11%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").Value = "'0.2666"
$ws.Range("E8").Value = "  +1.23%  "

$ws.Range("D9").Value = "'0.06605"
$ws.Range("E9").Value = "  +4.34%  "

$ws.Range("D10").Value = "'21.68"
$ws.Range("E10").Value = "  +6.46%  "

$ws.Range("D11").Value = "'0.07718"
$ws.Range("E11").Value = "  +0.81%  "

$ws.Range("D12").Value = "'4.621"
$ws.Range("E12").Value = "  +0.67%  "

$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.960.92"
$ws.Range("E13").Value = "  +4.45%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.705.99"
$ws.Range("E14").Value = "  +3.77%  "

$ws.Range("D15").Value = "'0.5839"
$ws.Range("E15").Value = "  +4.50%  "

$ws.Range("D16").Value = "0.0₅8307"
$ws.Range("E16").Value = "  +2.12%  "

$ws.Range("D17").Value = "'67.94"
$ws.Range("E17").Value = "  +4.24%  "

$ws.Range("D18").Value = "27.547.72"
$ws.Range("E18").Value = "  +5.74%  "

$ws.Range("D19").Value = "'220.05"
$ws.Range("E19").Value = "  +15.12%  "

$ws.Range("D20").Value = "'1.004"
$ws.Range("E20").Value = "  +0.09%  "

$ws.Range("D21").Value = "'4.728"
$ws.Range("E21").Value = "  +2.25%  "

$ws.Range("D22").Value = "'10.65"
$ws.Range("E22").Value = "  +1.35%  "

$ws.Range("D23").Value = "'6.083"
$ws.Range("E23").Value = "  +3.03%  "

$ws.Range("D24").Value = "'1.005"
$ws.Range("E24").Value = "  +0.17%  "

$ws.Range("D25").Value = "'148.12"
$ws.Range("E25").Value = "  +2.75%  "

$ws.Range("D26").Value = "'1.736"
$ws.Range("E26").Value = "  +14.62%  "

$ws.Range("D27").Value = "'0.1237"
$ws.Range("E27").Value = "  +4.36%  "

$ws.Range("D28").Value = "'7.400"
$ws.Range("E28").Value = "  +2.83%  "

$ws.Range("D29").Value = "'16.61"
$ws.Range("E29").Value = "  +4.59%  "

$ws.Range("D30").Value = "'0.05581"
$ws.Range("E30").Value = "  +2.68%  "

$ws.Range("E31").Value = "  +3.03%  "

$ws.Range("D32").Value = "'3.563"
$ws.Range("E32").Value = "  +3.67%  "

$ws.Range("D33").Value = "'3.447"
$ws.Range("E33").Value = "  +3.18%  "

$ws.Range("D34").Value = "'1.666"
$ws.Range("E34").Value = "  +7.15%  "

$ws.Range("E35").Value = "  +1.90%  "

$ws.Range("D36").Value = "'0.9614"
$ws.Range("E36").Value = "  +1.83%  "

$ws.Range("D37").Value = "'2.432"
$ws.Range("E37").Value = "  +0.25%  "

$ws.Range("D38").Value = "'0.5969"
$ws.Range("E38").Value = "  +5.92%  "

$ws.Range("D39").Value = "'0.01649"
$ws.Range("E39").Value = "  +4.46%  "

$ws.Range("D40").Value = "'5.924"
$ws.Range("E40").Value = "  +1.12%  "

$ws.Range("D41").Value = "'0.8544"
$ws.Range("E41").Value = "  +3.90%  "

$ws.Range("D42").Value = "1.052.45"
$ws.Range("E42").Value = "  +2.39%  "

$ws.Range("E43").Value = "  +0.14%  "

$ws.Range("D44").Value = "'101.33"
$ws.Range("E44").Value = "  +0.23%  "

$ws.Range("D45").Value = "1.867.01"
$ws.Range("E45").Value = "  +4.50%  "

$ws.Range("E46").Value = "  +3.93%  "

$ws.Range("D47").Value = "'59.19"
$ws.Range("E47").Value = "  +3.28%  "

$ws.Range("D48").Value = "'8.201"
$ws.Range("E48").Value = "  +3.13%  "

$ws.Range("D49").Value = "'0.4436"
$ws.Range("E49").Value = "  +2.42%  "

$ws.Range("D50").Value = "'1.001"
$ws.Range("E50").Value = "  -0.03%  "

$ws.Range("D51").Value = "'0.05245"
$ws.Range("E51").Value = "  +1.88%  "
